# Journal de travail - add "Développement du module Composal" entry and
# grow the work log table (Tableau1) to accommodate future rows, mirroring
# what Excel does when a user selects the table and drags its resize
# handle a number of rows further down, then types a new entry on the
# first of the freshly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")
$t = $ws.ListObjects.Item("Tableau1")

# --- 1. Fill in the new journal entry on the first currently-empty row ---
# Row 60 (A60) already exists as an empty table row; populate it.
$ws.Range("A60").Value = "2023-06-26"
$ws.Range("B60").Value = "Implémentation"
$ws.Range("D60").Value = "Développement du module Composal"

# --- 2. Grow the table by another 15 rows (new range A1:E77) ---
$t.Resize($ws.Range("A1:E77"))

# --- 3. Carry the "Date" column number formatting down into the freshly
#        appended rows, matching the rest of the column. ---
$ws.Range("A61:A75").NumberFormat = $ws.Range("A59").NumberFormat()

# --- 4. Move the active selection/view the way the author left it ---
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("D67").Select()
